$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Capture a style reference from a cell that the diff never touches (column C, a data-value style)
$origStyle = $ws.Range("C2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.63"
$ws.Range("D2").Style = $origStyle

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.05%"
$ws.Range("E2").Style = $origStyle

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.06"
$ws.Range("D3").Style = $origStyle

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.54%"
$ws.Range("E3").Style = $origStyle

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.508"
$ws.Range("D4").Style = $origStyle

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.97%"
$ws.Range("E4").Style = $origStyle

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08051"
$ws.Range("D5").Style = $origStyle

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.42%"
$ws.Range("E5").Style = $origStyle

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.106"
$ws.Range("D6").Style = $origStyle

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "9.80%"
$ws.Range("E6").Style = $origStyle

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9548"
$ws.Range("D7").Style = $origStyle

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.38%"
$ws.Range("E7").Style = $origStyle

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.62%"
$ws.Range("E8").Style = $origStyle

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1876"
$ws.Range("D9").Style = $origStyle

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.34%"
$ws.Range("E9").Style = $origStyle

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "10.17"
$ws.Range("D10").Style = $origStyle

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.75%"
$ws.Range("E10").Style = $origStyle

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09919"
$ws.Range("D11").Style = $origStyle

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.45%"
$ws.Range("E11").Style = $origStyle

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04730"
$ws.Range("D12").Style = $origStyle

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.25%"
$ws.Range("E12").Style = $origStyle

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.86%"
$ws.Range("E13").Style = $origStyle

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001267"
$ws.Range("D14").Style = $origStyle

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.19%"
$ws.Range("E14").Style = $origStyle

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04088"
$ws.Range("D15").Style = $origStyle

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.45%"
$ws.Range("E15").Style = $origStyle

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006051"
$ws.Range("D16").Style = $origStyle

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.43%"
$ws.Range("E16").Style = $origStyle

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.41%"
$ws.Range("E17").Style = $origStyle

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.95%"
$ws.Range("E18").Style = $origStyle

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.27%"
$ws.Range("E19").Style = $origStyle

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3307"
$ws.Range("D20").Style = $origStyle

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.50%"
$ws.Range("E20").Style = $origStyle

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.19%"
$ws.Range("E21").Style = $origStyle

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001309"
$ws.Range("D23").Style = $origStyle

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.22%"
$ws.Range("E23").Style = $origStyle

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004331"
$ws.Range("D24").Style = $origStyle

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.85%"
$ws.Range("E24").Style = $origStyle

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001252"
$ws.Range("D25").Style = $origStyle

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.26%"
$ws.Range("E25").Style = $origStyle

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-5.75%"
$ws.Range("E26").Style = $origStyle

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02617"
$ws.Range("D38").Style = $origStyle

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-2.55%"
$ws.Range("E38").Style = $origStyle

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05612"
$ws.Range("D39").Style = $origStyle

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.02%"
$ws.Range("E39").Style = $origStyle

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007606"
$ws.Range("D40").Style = $origStyle

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.64%"
$ws.Range("E40").Style = $origStyle

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1398"
$ws.Range("D41").Style = $origStyle

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.72%"
$ws.Range("E41").Style = $origStyle

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007376"
$ws.Range("D42").Style = $origStyle

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-11.23%"
$ws.Range("E42").Style = $origStyle

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001989"
$ws.Range("D43").Style = $origStyle

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.28%"
$ws.Range("E43").Style = $origStyle

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008855"
$ws.Range("D44").Style = $origStyle

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.00%"
$ws.Range("E44").Style = $origStyle

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007113"
$ws.Range("D45").Style = $origStyle

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.73%"
$ws.Range("E45").Style = $origStyle

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.33%"
$ws.Range("E46").Style = $origStyle

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005815"
$ws.Range("D47").Style = $origStyle

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "54.26%"
$ws.Range("E48").Style = $origStyle

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003365"
$ws.Range("D49").Style = $origStyle

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "6.68%"
$ws.Range("E49").Style = $origStyle

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("D50").Style = $origStyle

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.33%"
$ws.Range("E50").Style = $origStyle

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.33%"
$ws.Range("E51").Style = $origStyle

